$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3931.9092
$ws.Range("I33").Value = 750.7143
$ws.Range("K33").Value = 750.7143
$ws.Range("M33").Value = -521.7143
$ws.Range("H43").Value = 2742.6667
$ws.Range("I43").Value = 2506.45
$ws.Range("J43").Value = 3923.75
$ws.Range("K43").Value = 2506.45
$ws.Range("L43").Value = 3923.75
$ws.Range("M43").Value = -2437.45
$ws.Range("N43").Value = -4061.75
$ws.Range("H62").Value = 12373.385
$ws.Range("I62").Value = 10995.667
$ws.Range("K62").Value = 10995.667
$ws.Range("M62").Value = -10371.667
$ws.Range("H65").Value = 12373.385
$ws.Range("I65").Value = 10995.667
$ws.Range("K65").Value = 54978.335
$ws.Range("M65").Value = -51858.335
$ws.Range("H103").Value = 746.6875
$ws.Range("I103").Value = 743.55554
$ws.Range("K103").Value = 2230.66662
$ws.Range("M103").Value = -1644.66662
$ws.Range("H132").Value = 3938.8235
$ws.Range("I132").Value = 3797.3333
$ws.Range("K132").Value = 11391.9999
$ws.Range("M132").Value = -8861.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2615.8
$ws.Range("I2").Value = 860.55554
$ws.Range("K2").Value = 860.55554
$ws.Range("M2").Value = -747.55554
$ws.Range("H45").Value = 3706248
$ws.Range("I45").Value = 4547463.5
$ws.Range("J45").Value = 4899.6
$ws.Range("K45").Value = 4547463.5
$ws.Range("L45").Value = 4899.6
$ws.Range("M45").Value = -4547086.5
$ws.Range("N45").Value = -5653.6
$ws.Range("H61").Value = 7900.1313
$ws.Range("I61").Value = 7214.613
$ws.Range("J61").Value = 10936
$ws.Range("K61").Value = 7214.613
$ws.Range("L61").Value = 10936
$ws.Range("M61").Value = -7002.613
$ws.Range("N61").Value = -11360
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H116").Value = 2615.8
$ws.Range("I116").Value = 860.55554
$ws.Range("K116").Value = 860.55554
$ws.Range("M116").Value = 1433.44446
$ws.Range("H122").Value = 3105.1
$ws.Range("I122").Value = 2337.1667
$ws.Range("J122").Value = 10016.5
$ws.Range("K122").Value = 7011.500100000001
$ws.Range("L122").Value = 30049.5
$ws.Range("M122").Value = -4561.500100000001
$ws.Range("N122").Value = -34949.5
$ws.Range("H136").Value = 7900.1313
$ws.Range("I136").Value = 7214.613
$ws.Range("J136").Value = 10936
$ws.Range("K136").Value = 21643.839
$ws.Range("L136").Value = 32808
$ws.Range("M136").Value = -19093.839
$ws.Range("N136").Value = -37908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2615.8
$ws.Range("I3").Value = 860.55554
$ws.Range("K3").Value = 860.55554
$ws.Range("M3").Value = -746.55554
$ws.Range("H33").Value = 121500
$ws.Range("I33").Value = 87500
$ws.Range("J33").Value = 189500
$ws.Range("K33").Value = 87500
$ws.Range("L33").Value = 189500
$ws.Range("M33").Value = -87164
$ws.Range("N33").Value = -190172
$ws.Range("H94").Value = 1283.7333
$ws.Range("I94").Value = 604.8333
$ws.Range("J94").Value = 3999.3333
$ws.Range("K94").Value = 604.8333
$ws.Range("L94").Value = 3999.3333
$ws.Range("M94").Value = -153.8333
$ws.Range("N94").Value = -4901.3333
$ws.Range("H105").Value = 757988.2
$ws.Range("I105").Value = 2671829.8
$ws.Range("J105").Value = 2524.4211
$ws.Range("K105").Value = 2671829.8
$ws.Range("L105").Value = 2524.4211
$ws.Range("M105").Value = -2670082.8
$ws.Range("N105").Value = -6018.4211
$ws.Range("H134").Value = 6726.9443
$ws.Range("I134").Value = 4780.0967
$ws.Range("J134").Value = 18797.4
$ws.Range("K134").Value = 14340.2901
$ws.Range("L134").Value = 56392.2
$ws.Range("M134").Value = -11805.2901
$ws.Range("N134").Value = -61462.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 13199.667
$ws.Range("J6").Value = 38000
$ws.Range("L6").Value = 38000
$ws.Range("N6").Value = -38226
$ws.Range("H19").Value = 30668.666
$ws.Range("H24").Value = 30668.666
$ws.Range("H32").Value = 7224.75
$ws.Range("I32").Value = 7449.5
$ws.Range("K32").Value = 7449.5
$ws.Range("M32").Value = -7133.5
$ws.Range("H50").Value = 53999.5
$ws.Range("J50").Value = 53999.5
$ws.Range("L50").Value = 53999.5
$ws.Range("N50").Value = -55249.5
$ws.Range("H58").Value = 5297.8184
$ws.Range("I58").Value = 4744.857
$ws.Range("K58").Value = 4744.857
$ws.Range("M58").Value = -4541.857
$ws.Range("H59").Value = 79261.375
$ws.Range("I59").Value = 20052
$ws.Range("J59").Value = 98997.836
$ws.Range("K59").Value = 20052
$ws.Range("L59").Value = 98997.836
$ws.Range("M59").Value = -18907
$ws.Range("N59").Value = -101287.836
$ws.Range("H95").Value = 35465.555
$ws.Range("J95").Value = 35465.555
$ws.Range("L95").Value = 35465.555
$ws.Range("N95").Value = -40957.555
$ws.Range("H132").Value = 5780.1777
$ws.Range("I132").Value = 4301.278
$ws.Range("K132").Value = 12903.834
$ws.Range("M132").Value = -10373.834
$ws.Range("H134").Value = 4242.3335
$ws.Range("I134").Value = 3342.9473
$ws.Range("J134").Value = 7660
$ws.Range("K134").Value = 10028.8419
$ws.Range("L134").Value = 22980
$ws.Range("M134").Value = -7493.841899999999
$ws.Range("N134").Value = -28050
$ws.Range("H136").Value = 5297.8184
$ws.Range("I136").Value = 4744.857
$ws.Range("K136").Value = 14234.571
$ws.Range("M136").Value = -11684.571
$ws.Range("H141").Value = 200447
$ws.Range("J141").Value = 245404.75
$ws.Range("L141").Value = 245404.75
$ws.Range("N141").Value = -255764.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 51.916668
$ws.Range("I2").Value = 66.25
$ws.Range("J2").Value = 23.25
$ws.Range("K2").Value = 397.5
$ws.Range("L2").Value = 139.5
$ws.Range("M2").Value = -284.5
$ws.Range("N2").Value = -365.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6214.4
$ws.Range("I122").Value = 3984.0667
$ws.Range("J122").Value = 9559.9
$ws.Range("K122").Value = 11952.2001
$ws.Range("L122").Value = 28679.7
$ws.Range("M122").Value = -9502.2001
$ws.Range("N122").Value = -33579.7
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4495.857
$ws.Range("I7").Value = 4495.857
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4495.857
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4383.857
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 3273
$ws.Range("I22").Value = 3670.6667
$ws.Range("J22").Value = 2761.7144
$ws.Range("K22").Value = 3670.6667
$ws.Range("L22").Value = 2761.7144
$ws.Range("M22").Value = -3375.6667
$ws.Range("N22").Value = -3351.7144
$ws.Range("H27").Value = 3273
$ws.Range("I27").Value = 3670.6667
$ws.Range("J27").Value = 2761.7144
$ws.Range("K27").Value = 3670.6667
$ws.Range("L27").Value = 2761.7144
$ws.Range("M27").Value = -3563.6667
$ws.Range("N27").Value = -2975.7144
$ws.Range("H122").Value = 3694.6155
$ws.Range("I122").Value = 1853.5
$ws.Range("K122").Value = 5560.5
$ws.Range("M122").Value = -3110.5
$ws.Range("H126").Value = 4495.857
$ws.Range("I126").Value = 4495.857
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13487.571
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11017.571
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 10086.407
$ws.Range("I136").Value = 8641.076999999999
$ws.Range("J136").Value = 11428.5
$ws.Range("K136").Value = 25923.231
$ws.Range("L136").Value = 34285.5
$ws.Range("M136").Value = -23373.231
$ws.Range("N136").Value = -39385.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H122").Value = 2253.9
$ws.Range("I122").Value = 2210.9312
$ws.Range("K122").Value = 6632.7936
$ws.Range("M122").Value = -4182.7936
$ws.Range("H126").Value = 3420.8462
$ws.Range("I126").Value = 1433.375
$ws.Range("J126").Value = 6600.8
$ws.Range("K126").Value = 4300.125
$ws.Range("L126").Value = 19802.4
$ws.Range("M126").Value = -1830.125
$ws.Range("N126").Value = -24742.4
